# Workbook has two sheets: "Planilha1" (a small table with a header row
# plus 6 data rows, A1:C7) and "Planilha2" (a single, empty, underlined
# cell A1).
#
# This edit:
#   1. Clears all the text out of the Planilha1 table (A1:C7) while
#      leaving each cell's formatting (fills/borders/fonts) in place.
#   2. Adds an underline on top of the existing header font (row 1,
#      bold) and the existing body font (rows 2-7, italic).
#   3. Makes Planilha1 the active sheet, with the whole table (A1:C7)
#      selected and the view scrolled back to the top-left.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Planilha1")

# 1. Wipe the cell text but keep every cell's existing style/formatting.
$ws1.Range("A1:C7").ClearContents() | Out-Null

# 2. Underline the header row and the body rows (each keeps its own
#    bold / italic styling, just with underline added on top).
$ws1.Range("A1:C1").Font.Underline = $true
$ws1.Range("A2:C7").Font.Underline = $true

# 3. Activate Planilha1, scroll back to the top-left, and select the
#    full table range (this also drops Planilha2's tabSelected flag
#    and the workbook's explicit activeTab, since Planilha1 - index 0 -
#    becomes the active / selected tab).
$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("A1:C7").Select() | Out-Null
